$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 (a full clinical-data record) into a new row 3,
# preserving its styling/number formats, then give the new row its
# own subject/sample identifier in column A.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "test_subject__test_sample"

# Move the active selection to reflect the newly added row.
$ws.Range("B4").Select() | Out-Null
